# Se realizan cambios para sanity semilla 8
#
# - "Semilla 3" sheet becomes "Semilla 9" with refreshed server/test data.
# - A brand-new "Semilla 8" sheet is added (copied from the original
#   "Semilla 3" sheet so it keeps the same structure/formatting/hyperlinks)
#   with its own refreshed server/test data.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)      # "Semilla 4" - untouched
$ws2 = $wb.Worksheets.Item(2)      # "Semilla 3" -> "Semilla 9"

# ---------------------------------------------------------------------
# 1) Create the new "Semilla 8" sheet as a full copy of the current
#    "Semilla 3" sheet (BEFORE we touch its data), placed right after it.
#    This preserves formatting, column widths, and the 5 separate
#    hyperlinks (on A2,B2,C2,D2,E2) that the new sheet needs.
# ---------------------------------------------------------------------
$ws2.Copy($null, $ws2)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "Semilla 8"

# Update the hyperlink addresses on "Semilla 8" for the cells whose text
# changes (A2, B2, C2) so the link target matches the new visible text.
foreach ($h in $ws3.Hyperlinks) {
    $col = $h.Range.Column
    if ($col -eq 1) { $h.Address = "http://10.69.60.121:8180/tigo-pos-web/index.jsp" }
    elseif ($col -eq 2) { $h.Address = "http://10.69.60.123:8280/portal/login?initialURI=%2Fportal%2FCRMPortal" }
    elseif ($col -eq 3) { $h.Address = "http://10.69.60.121:8180/tigo-pos-web/wap/windex.wml" }
}

# Refresh the visible data on "Semilla 8".
$ws3.Range("A2").Value2 = "http://10.69.60.121:8180/tigo-pos-web/index.jsp"
$ws3.Range("B2").Value2 = "http://10.69.60.123:8280/portal/login?initialURI=%2Fportal%2FCRMPortal"
$ws3.Range("C2").Value2 = "http://10.69.60.121:8180/tigo-pos-web/wap/windex.wml"

$ws3.Range("A4").Value2 = "10.69.60.115"

$ws3.Range("A5").Value2 = "10.69.60.113"
$ws3.Range("B5").Value2 = "DESEPOS"

$ws3.Range("A6").Value2 = "10.69.60.113"
$ws3.Range("B6").Value2 = "DEV10G"

$ws3.Range("B7").Value2 = "siebel04"

$ws3.Range("B9").Value2  = "212238780"
$ws3.Range("B10").Value2 = "250453824"
$ws3.Range("B11").Value2 = "664766732"
$ws3.Range("B12").Value2 = "220291650"
$ws3.Range("B13").Value2 = "125782671"

# New sheet is not the selected tab; its own selection cursor is B16.
$ws3.Range("B16").Select()

# ---------------------------------------------------------------------
# 2) Rename "Semilla 3" to "Semilla 9" and refresh its data in place.
# ---------------------------------------------------------------------
$ws2.Name = "Semilla 9"

# Update hyperlink addresses for A2/B2/C2 before collapsing them down to
# the single surviving hyperlink on B2.
foreach ($h in $ws2.Hyperlinks) {
    $col = $h.Range.Column
    if ($col -eq 1) { $h.Address = "http://10.69.60.137:8180/tigo-pos-web/index.jsp" }
    elseif ($col -eq 2) { $h.Address = "http://10.69.60.140:8280/portal/login?initialURI=%2Fportal%2FCRMPortal%2FVenta" }
    elseif ($col -eq 3) { $h.Address = "http://10.69.60.137:8180/tigo-pos-web/wap/windex.wml" }
}

# "Semilla 9" ends up with only one hyperlink left: B2.
$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("B2"), "http://10.69.60.140:8280/portal/login?initialURI=%2Fportal%2FCRMPortal%2FVenta")

$ws2.Range("A2").Value2 = "http://10.69.60.137:8180/tigo-pos-web/index.jsp"
$ws2.Range("B2").Value2 = "http://10.69.60.140:8280/portal/login?initialURI=%2Fportal%2FCRMPortal%2FVenta"
$ws2.Range("C2").Value2 = "http://10.69.60.137:8180/tigo-pos-web/wap/windex.wml"

$ws2.Range("A4").Value2 = "10.69.60.132"
$ws2.Range("B4").Value2 = "dev11g"

$ws2.Range("A5").Value2 = "10.69.60.130"
$ws2.Range("B5").Value2 = "DESEPOS"

$ws2.Range("A6").Value2 = "10.69.60.130"

$ws2.Range("A7").Value2 = "10.65.32.74"
$ws2.Range("B7").Value2 = "SIEBEL05"
$ws2.Range("C7").Value2 = "SIEBEL"

$ws2.Range("B9").Value2  = "553312726"
$ws2.Range("B10").Value2 = "15377510"
$ws2.Range("B11").Value2 = "410614432"
$ws2.Range("B12").Value2 = "22368093"
$ws2.Range("B13").Value2 = "1050388676"

# "Semilla 9" stays the active/selected tab, with the cursor now on B3.
$ws2.Activate()
$ws2.Range("B3").Select()
